$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.978.25"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  +3.42%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D12").Value = "1.777.84"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "1.556.54"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").Value = "26.962.81"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.78"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.18"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.95"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "1.430.11"
$ws.Range("E33").Value = "  +5.28%  "
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("E35").Value = "  +3.87%  "
$ws.Range("E36").Value = "  +2.10%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("E43").Value = "  +3.87%  "
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.59"
$ws.Range("D45").ClearFormats()
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("D47").Value = "1.691.50"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.75"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("E50").Value = "  +3.14%  "
$ws.Range("E51").Value = "  +1.30%  "
